# The commit inserts one new data row (a new weekly "Ajo" price record for
# "Macroferia Regional de Talca") right before the existing row that used to
# be row 279. All subsequent rows shift down by one (old row N -> new row
# N+1), and the sheet's used-range dimension grows from A1:R340 to A1:R341.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 279; everything from the old row 279
# downwards (through row 340) is pushed down to rows 280-341.
$ws.Rows("279").Insert()

# Populate the newly inserted row 279 with the new record's data.
$ws.Range("A279").Value = 5
$ws.Range("B279").Value = "Macroferia Regional de Talca"
$ws.Range("C279").Value = "Maule"
$ws.Range("D279").Value = 44785
$ws.Range("E279").Value = 7
$ws.Range("F279").Value = 100112003
$ws.Range("G279").Value = "Ajo"
$ws.Range("H279").Value = "Chino"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 300
$ws.Range("K279").Value = 28000
$ws.Range("L279").Value = 28000
$ws.Range("M279").Value = 28000
$ws.Range("N279").Value = "`$/malla 10 kilos"
$ws.Range("O279").Value = "China"
$ws.Range("P279").Value = 2800
$ws.Range("Q279").Value = 10
$ws.Range("R279").Value = "Hortaliza"
